# Correct typos & update offloading fig
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct typos: add trailing period to the "Only ..." comment strings.
$ws.Range("B5").Value  = "Only training=False."
$ws.Range("B7").Value  = "Only along channel axis."
$ws.Range("B10").Value = "Only constant multiplication."
$ws.Range("B14").Value = "Only inference mode."
$ws.Range("B15").Value = "Only channel-wise flatten and before fully connected layer or Conv w/ 1x1 kernel."
$ws.Range("B16").Value = "Only alpha=0, max_value=None, threshold=0."

# Update selection to reflect the cell near the updated offloading figure.
$ws.Range("B18").Select()
